$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# 1) Widen the two "p-value" grid columns (4 and 7) from 1132 -> 1181 twips
$t.Columns.Item(4).Width = 1181 / 20.0
$t.Columns.Item(7).Width = 1181 / 20.0

# 2) Bump the second header row's height from 614 -> 615 twips
$headerRow = $t.Rows.Item(2)
$headerRow.Height = 615 / 20.0

# 3) Bold the main label in every cell of that header row, leaving the
#    superscript footnote markers ("1") un-bolded.
for ($c = 1; $c -le $headerRow.Cells.Count; $c++) {
    $cell = $headerRow.Cells.Item($c)
    $start = $cell.Range.Start
    $end = $cell.Range.End - 1   # drop the trailing end-of-cell mark

    $boldEnd = $end
    for ($i = $start; $i -lt $end; $i++) {
        $ch = $d.Range($i, $i + 1)
        if ($ch.Font.Superscript) {
            $boldEnd = $i
            break
        }
    }

    if ($boldEnd -gt $start) {
        $d.Range($start, $boldEnd).Font.Bold = $true
    }
}
